$d = $word.ActiveDocument

$d.Content.Find.Execute("948×9=8532", $true, $false, $false, $false, $false, $true, 1, $false, "119×9=1071", 2) | Out-Null
$d.Content.Find.Execute("433×9=3897", $true, $false, $false, $false, $false, $true, 1, $false, "480×7=3360", 2) | Out-Null
$d.Content.Find.Execute("171×2=342", $true, $false, $false, $false, $false, $true, 1, $false, "148×8=1184", 2) | Out-Null
$d.Content.Find.Execute("419×5=2095", $true, $false, $false, $false, $false, $true, 1, $false, "491×3=1473", 2) | Out-Null
$d.Content.Find.Execute("927×9=8343", $true, $false, $false, $false, $false, $true, 1, $false, "954×4=3816", 2) | Out-Null
$d.Content.Find.Execute("601×3=1803", $true, $false, $false, $false, $false, $true, 1, $false, "272×9=2448", 2) | Out-Null
$d.Content.Find.Execute("389×7=2723", $true, $false, $false, $false, $false, $true, 1, $false, "154×7=1078", 2) | Out-Null
$d.Content.Find.Execute("921×7=6447", $true, $false, $false, $false, $false, $true, 1, $false, "201×7=1407", 2) | Out-Null
$d.Content.Find.Execute("463×8=3704", $true, $false, $false, $false, $false, $true, 1, $false, "403×5=2015", 2) | Out-Null
$d.Content.Find.Execute("119×4=476", $true, $false, $false, $false, $false, $true, 1, $false, "952×8=7616", 2) | Out-Null
$d.Content.Find.Execute("957×4=3828", $true, $false, $false, $false, $false, $true, 1, $false, "752×8=6016", 2) | Out-Null
$d.Content.Find.Execute("375×8=3000", $true, $false, $false, $false, $false, $true, 1, $false, "997×7=6979", 2) | Out-Null
$d.Content.Find.Execute("248×8=1984", $true, $false, $false, $false, $false, $true, 1, $false, "648×3=1944", 2) | Out-Null
$d.Content.Find.Execute("536×2=1072", $true, $false, $false, $false, $false, $true, 1, $false, "587×8=4696", 2) | Out-Null
$d.Content.Find.Execute("351×3=1053", $true, $false, $false, $false, $false, $true, 1, $false, "388×2=776", 2) | Out-Null
$d.Content.Find.Execute("208×3=624", $true, $false, $false, $false, $false, $true, 1, $false, "597×7=4179", 2) | Out-Null
$d.Content.Find.Execute("130×4=520", $true, $false, $false, $false, $false, $true, 1, $false, "354×6=2124", 2) | Out-Null
$d.Content.Find.Execute("708×6=4248", $true, $false, $false, $false, $false, $true, 1, $false, "188×3=564", 2) | Out-Null
$d.Content.Find.Execute("258×5=1290", $true, $false, $false, $false, $false, $true, 1, $false, "205×6=1230", 2) | Out-Null
$d.Content.Find.Execute("971×5=4855", $true, $false, $false, $false, $false, $true, 1, $false, "900×4=3600", 2) | Out-Null
$d.Content.Find.Execute("710×3=2130", $true, $false, $false, $false, $false, $true, 1, $false, "495×8=3960", 2) | Out-Null
$d.Content.Find.Execute("868×5=4340", $true, $false, $false, $false, $false, $true, 1, $false, "880×9=7920", 2) | Out-Null
$d.Content.Find.Execute("331×6=1986", $true, $false, $false, $false, $false, $true, 1, $false, "223×3=669", 2) | Out-Null
$d.Content.Find.Execute("406×2=812", $true, $false, $false, $false, $false, $true, 1, $false, "155×3=465", 2) | Out-Null
$d.Content.Find.Execute("835×5=4175", $true, $false, $false, $false, $false, $true, 1, $false, "110×9=990", 2) | Out-Null
